# Adds three new programs (row 21-23) to the programs table:
#   - Job Creation Schemes ("Arbeitsbeschaffungsmaßnahmen")
#   - One Euro Jobs ("Ein Euro Jobs")
#   - Subsidized Job Opportunities ("Arbeitsgelegenheiten")
# based on Hohmeyer & Wolff (2010).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Job Creation Schemes ------------------------------------------------
$ws.Cells.Item(21, 2).Value = "Job Creation Schemes"
$ws.Cells.Item(22, 1).Value = "oneEuroJobs"
$ws.Cells.Item(21, 1).Value = "jobCreationSchemes"
$ws.Cells.Item(22, 6).Value = 'One Euro Jobs "Ein Euro Jobs" are jobs which are assigned to long term unemployed to help them reintegrate into the labor market. These jobs are meant to be unpayed. However, the participants receive a compenstation of at leat one euro per hour worked.'
$ws.Cells.Item(22, 2).Value = "One Euro Jobs"
$ws.Cells.Item(21, 6).Value = 'Job Creation Schemes "Arbeitsbeschaffungsmaßnahmen" are comprised of subsizied jobs aimed at  unemployed who would not be able to find a job otherwise. This measure is mostly employed in local labor markets with excess labor supply. The duration of the subsidy (between 900 and 1300€ payed to the employer) is usually limited to 12 months.'
$ws.Cells.Item(23, 1).Value = "subsidizedJobOpportunities"
$ws.Cells.Item(23, 2).Value = "Subsidized Job Opportunities"
$ws.Cells.Item(23, 6).Value = 'Similar to  Job Creation Schemes, Subsidized Job Opportunities "Arbeitsgelegenheiten" are jobs where the employer receives a subsidy. The focus of subsidized job opportunities lies on reintegration. The amount of subsidy payed to the employer is not fixed and the types of jobs that can be subsizied  is more broad. '

# --- year ------------------------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 2005
$ws.Cells.Item(22, 3).Value = 2005
$ws.Cells.Item(23, 3).Value = 2005

# --- category ----------------------------------------------------------------------
$ws.Cells.Item(21, 4).Value = "Active Labor Market Policy"
$ws.Cells.Item(22, 4).Value = "Active Labor Market Policy"
$ws.Cells.Item(23, 4).Value = "Active Labor Market Policy"

# --- row heights (wrapped long-form description text) -----------------------------
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(23).RowHeight = 105

# --- move the active selection to the new description cell ------------------------
$ws.Range("E21").Select() | Out-Null
